$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Price")

# --- Header row ---
$ws.Cells.Item(1,2).Value = "fuel_type_id"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Price"
$ws.Cells.Item(1,5).Value = "Effective_date"

# --- Data rows: Price_id | fuel_type_id | Name | Price | Effective_date ---
# Row 2
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = "Diesel"
$ws.Cells.Item(2,4).Value = 60
$ws.Cells.Item(2,5).Value = "date"

# Row 3
$ws.Cells.Item(3,2).Value = 2
$ws.Cells.Item(3,3).Value = "Premium"
$ws.Cells.Item(3,4).Value = 50
$ws.Cells.Item(3,5).Value = "date"

# Row 4
$ws.Cells.Item(4,2).Value = 3
$ws.Cells.Item(4,3).Value = "Unleaded"
$ws.Cells.Item(4,4).Value = 70
$ws.Cells.Item(4,5).Value = "date"

# Row 5
$ws.Cells.Item(5,2).Value = 1
$ws.Cells.Item(5,3).Value = "Diesel100"
$ws.Cells.Item(5,4).Value = 60
$ws.Cells.Item(5,5).Value = "date"

# Row 6
$ws.Cells.Item(6,2).Value = 2
$ws.Cells.Item(6,3).Value = "Premium100"
$ws.Cells.Item(6,4).Value = 60
$ws.Cells.Item(6,5).Value = "date"

# Row 7
$ws.Cells.Item(7,2).Value = 3
$ws.Cells.Item(7,3).Value = "Unleaded100"
$ws.Cells.Item(7,4).Value = 40
$ws.Cells.Item(7,5).Value = "date"

# Row 8
$ws.Cells.Item(8,2).Value = 1
$ws.Cells.Item(8,3).Value = "Diesel100"
$ws.Cells.Item(8,4).Value = 60
$ws.Cells.Item(8,5).Value = "date"

# Row 9
$ws.Cells.Item(9,2).Value = 2
$ws.Cells.Item(9,3).Value = "Premium"
$ws.Cells.Item(9,4).Value = 60
$ws.Cells.Item(9,5).Value = "date"

# --- New column E width ---
$ws.Columns("E").ColumnWidth = 12.6745

# --- Make "Price" the active sheet/tab and set the selection ---
$ws.Activate()
$ws.Range("D11").Select() | Out-Null
